# Update NATMI ligand-receptor pair TPM statistics (Lamc2-Itga6) with
# recomputed values after the underlying TPM data was refreshed.
# Values are set via [double]"..." casts so that high-precision decimal
# and scientific-notation literals round-trip exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"0.7955796666666667"
$ws.Range("H2").Value = [double]"2.386739"
$ws.Range("I2").Value = [double]"0.1186174580157865"
$ws.Range("J2").Value = [double]"0.1186174580157865"
$ws.Range("M2").Value = [double]"227.11144"
$ws.Range("N2").Value = [double]"681.33432"
$ws.Range("O2").Value = [double]"0.8625743548356182"
$ws.Range("P2").Value = [double]"0.8625743548356182"
$ws.Range("Q2").Value = [double]"180.6852437313867"
$ws.Range("R2").Value = [double]"1626.16719358248"
$ws.Range("S2").Value = [double]"0.1023163773202081"
$ws.Range("T2").Value = [double]"0.1023163773202081"
$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"0.7955796666666667"
$ws.Range("H3").Value = [double]"2.386739"
$ws.Range("I3").Value = [double]"0.1186174580157865"
$ws.Range("J3").Value = [double]"0.1186174580157865"
$ws.Range("O3").Value = [double]"0.001598666154760757"
$ws.Range("P3").Value = [double]"0.001598666154760757"
$ws.Range("Q3").Value = [double]"0.3348759236797778"
$ws.Range("R3").Value = [double]"3.013883313118"
$ws.Range("S3").Value = [double]"0.000189629715493593"
$ws.Range("T3").Value = [double]"0.0001896297154935929"
$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"0.7955796666666667"
$ws.Range("H4").Value = [double]"2.386739"
$ws.Range("I4").Value = [double]"0.1186174580157865"
$ws.Range("J4").Value = [double]"0.1186174580157865"
$ws.Range("M4").Value = [double]"3.233093"
$ws.Range("N4").Value = [double]"9.699279000000001"
$ws.Range("O4").Value = [double]"0.01227935989749593"
$ws.Range("P4").Value = [double]"0.01227935989749593"
$ws.Range("Q4").Value = [double]"2.572183051242333"
$ws.Range("R4").Value = [double]"23.149647461181"
$ws.Range("S4").Value = [double]"0.001456546457101957"
$ws.Range("T4").Value = [double]"0.001456546457101956"
$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"0.7955796666666667"
$ws.Range("H5").Value = [double]"2.386739"
$ws.Range("I5").Value = [double]"0.1186174580157865"
$ws.Range("J5").Value = [double]"0.1186174580157865"
$ws.Range("M5").Value = [double]"32.52945966666667"
$ws.Range("N5").Value = [double]"97.588379"
$ws.Range("O5").Value = [double]"0.1235476191121251"
$ws.Range("P5").Value = [double]"0.1235476191121251"
$ws.Range("Q5").Value = [double]"25.87977667845345"
$ws.Range("R5").Value = [double]"232.917990106081"
$ws.Range("S5").Value = [double]"0.01465490452298289"
$ws.Range("T5").Value = [double]"0.01465490452298289"
$ws.Range("I6").Value = [double]"0.6312226244877757"
$ws.Range("J6").Value = [double]"0.6312226244877758"
$ws.Range("M6").Value = [double]"227.11144"
$ws.Range("N6").Value = [double]"681.33432"
$ws.Range("O6").Value = [double]"0.8625743548356182"
$ws.Range("P6").Value = [double]"0.8625743548356182"
$ws.Range("Q6").Value = [double]"961.5162528534401"
$ws.Range("R6").Value = [double]"8653.646275680961"
$ws.Range("S6").Value = [double]"0.5444764480751888"
$ws.Range("T6").Value = [double]"0.5444764480751889"
$ws.Range("I7").Value = [double]"0.6312226244877757"
$ws.Range("J7").Value = [double]"0.6312226244877758"
$ws.Range("O7").Value = [double]"0.001598666154760757"
$ws.Range("P7").Value = [double]"0.001598666154760757"
$ws.Range("S7").Value = [double]"0.001009114245887865"
$ws.Range("T7").Value = [double]"0.001009114245887866"
$ws.Range("I8").Value = [double]"0.6312226244877757"
$ws.Range("J8").Value = [double]"0.6312226244877758"
$ws.Range("M8").Value = [double]"3.233093"
$ws.Range("N8").Value = [double]"9.699279000000001"
$ws.Range("O8").Value = [double]"0.01227935989749593"
$ws.Range("P8").Value = [double]"0.01227935989749593"
$ws.Range("Q8").Value = [double]"13.687868239868"
$ws.Range("R8").Value = [double]"123.190814158812"
$ws.Range("S8").Value = [double]"0.007751009781527326"
$ws.Range("T8").Value = [double]"0.007751009781527326"
$ws.Range("I9").Value = [double]"0.6312226244877757"
$ws.Range("J9").Value = [double]"0.6312226244877758"
$ws.Range("M9").Value = [double]"32.52945966666667"
$ws.Range("N9").Value = [double]"97.588379"
$ws.Range("O9").Value = [double]"0.1235476191121251"
$ws.Range("P9").Value = [double]"0.1235476191121251"
$ws.Range("Q9").Value = [double]"137.7191926837347"
$ws.Range("R9").Value = [double]"1239.472734153612"
$ws.Range("S9").Value = [double]"0.0779860523851717"
$ws.Range("T9").Value = [double]"0.07798605238517171"
$ws.Range("G10").Value = [double]"1.666370333333333"
$ws.Range("H10").Value = [double]"4.999111"
$ws.Range("I10").Value = [double]"0.248448548064433"
$ws.Range("J10").Value = [double]"0.248448548064433"
$ws.Range("M10").Value = [double]"227.11144"
$ws.Range("N10").Value = [double]"681.33432"
$ws.Range("O10").Value = [double]"0.8625743548356182"
$ws.Range("P10").Value = [double]"0.8625743548356182"
$ws.Range("Q10").Value = [double]"378.4517659766133"
$ws.Range("R10").Value = [double]"3406.06589378952"
$ws.Range("S10").Value = [double]"0.2143053460565243"
$ws.Range("T10").Value = [double]"0.2143053460565243"
$ws.Range("G11").Value = [double]"1.666370333333333"
$ws.Range("H11").Value = [double]"4.999111"
$ws.Range("I11").Value = [double]"0.248448548064433"
$ws.Range("J11").Value = [double]"0.248448548064433"
$ws.Range("O11").Value = [double]"0.001598666154760757"
$ws.Range("P11").Value = [double]"0.001598666154760757"
$ws.Range("Q11").Value = [double]"0.7014097116202223"
$ws.Range("R11").Value = [double]"6.312687404582001"
$ws.Range("S11").Value = [double]"0.0003971862849900601"
$ws.Range("T11").Value = [double]"0.0003971862849900601"
$ws.Range("G12").Value = [double]"1.666370333333333"
$ws.Range("H12").Value = [double]"4.999111"
$ws.Range("I12").Value = [double]"0.248448548064433"
$ws.Range("J12").Value = [double]"0.248448548064433"
$ws.Range("M12").Value = [double]"3.233093"
$ws.Range("N12").Value = [double]"9.699279000000001"
$ws.Range("O12").Value = [double]"0.01227935989749593"
$ws.Range("P12").Value = [double]"0.01227935989749593"
$ws.Range("Q12").Value = [double]"5.387530260107667"
$ws.Range("R12").Value = [double]"48.487772340969"
$ws.Range("S12").Value = [double]"0.003050789137693488"
$ws.Range("T12").Value = [double]"0.003050789137693488"
$ws.Range("G13").Value = [double]"1.666370333333333"
$ws.Range("H13").Value = [double]"4.999111"
$ws.Range("I13").Value = [double]"0.248448548064433"
$ws.Range("J13").Value = [double]"0.248448548064433"
$ws.Range("M13").Value = [double]"32.52945966666667"
$ws.Range("N13").Value = [double]"97.588379"
$ws.Range("O13").Value = [double]"0.1235476191121251"
$ws.Range("P13").Value = [double]"0.1235476191121251"
$ws.Range("Q13").Value = [double]"54.20612654789655"
$ws.Range("R13").Value = [double]"487.855138931069"
$ws.Range("S13").Value = [double]"0.03069522658522508"
$ws.Range("T13").Value = [double]"0.03069522658522508"
$ws.Range("E14").Value = [double]"2"
$ws.Range("F14").Value = [double]"0.6666666666666666"
$ws.Range("G14").Value = [double]"0.01147833333333333"
$ws.Range("H14").Value = [double]"0.034435"
$ws.Range("I14").Value = [double]"0.00171136943200476"
$ws.Range("J14").Value = [double]"0.00171136943200476"
$ws.Range("M14").Value = [double]"227.11144"
$ws.Range("N14").Value = [double]"681.33432"
$ws.Range("O14").Value = [double]"0.8625743548356182"
$ws.Range("P14").Value = [double]"0.8625743548356182"
$ws.Range("Q14").Value = [double]"2.606860812133334"
$ws.Range("R14").Value = [double]"23.4617473092"
$ws.Range("S14").Value = [double]"0.001476183383696904"
$ws.Range("T14").Value = [double]"0.001476183383696904"
$ws.Range("E15").Value = [double]"2"
$ws.Range("F15").Value = [double]"0.6666666666666666"
$ws.Range("G15").Value = [double]"0.01147833333333333"
$ws.Range("H15").Value = [double]"0.034435"
$ws.Range("I15").Value = [double]"0.00171136943200476"
$ws.Range("J15").Value = [double]"0.00171136943200476"
$ws.Range("O15").Value = [double]"0.001598666154760757"
$ws.Range("P15").Value = [double]"0.001598666154760757"
$ws.Range("Q15").Value = [double]"0.004831467718888889"
$ws.Range("R15").Value = [double]"0.04348320947000001"
$ws.Range("S15").Value = [double]"2.735908389238151E-06"
$ws.Range("T15").Value = [double]"2.73590838923815E-06"
$ws.Range("E16").Value = [double]"2"
$ws.Range("F16").Value = [double]"0.6666666666666666"
$ws.Range("G16").Value = [double]"0.01147833333333333"
$ws.Range("H16").Value = [double]"0.034435"
$ws.Range("I16").Value = [double]"0.00171136943200476"
$ws.Range("J16").Value = [double]"0.00171136943200476"
$ws.Range("M16").Value = [double]"3.233093"
$ws.Range("N16").Value = [double]"9.699279000000001"
$ws.Range("O16").Value = [double]"0.01227935989749593"
$ws.Range("P16").Value = [double]"0.01227935989749593"
$ws.Range("Q16").Value = [double]"0.03711051915166667"
$ws.Range("R16").Value = [double]"0.333994672365"
$ws.Range("S16").Value = [double]"2.101452117315964E-05"
$ws.Range("T16").Value = [double]"2.101452117315964E-05"
$ws.Range("E17").Value = [double]"2"
$ws.Range("F17").Value = [double]"0.6666666666666666"
$ws.Range("G17").Value = [double]"0.01147833333333333"
$ws.Range("H17").Value = [double]"0.034435"
$ws.Range("I17").Value = [double]"0.00171136943200476"
$ws.Range("J17").Value = [double]"0.00171136943200476"
$ws.Range("M17").Value = [double]"32.52945966666667"
$ws.Range("N17").Value = [double]"97.588379"
$ws.Range("O17").Value = [double]"0.1235476191121251"
$ws.Range("P17").Value = [double]"0.1235476191121251"
$ws.Range("Q17").Value = [double]"0.3733839812072222"
$ws.Range("R17").Value = [double]"3.360455830865"
$ws.Range("S17").Value = [double]"0.0002114356187454581"
$ws.Range("T17").Value = [double]"0.0002114356187454581"
